# Fix bugs during test, update test progress (finished basic data types test),
# add field logging (UNTESTED)
#
# The data table on Sheet1 gains 8 new config rows (client/server request &
# response logging toggles, each with a matching "print-fields" toggle) that
# are sorted into place right after the existing "core.test.*" rows, i.e.
# inserted at row 13, pushing the previous rows 13-45 down to 21-53.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# ---------------------------------------------------------------------------
# 1. Make room for the 8 new rows at the top of the block (old row 13..45
#    becomes 21..53), carrying the existing formatting down with them.
# ---------------------------------------------------------------------------
$ws.Rows("13:20").Insert()

# Newly inserted rows don't inherit the sheet's standard row height
# automatically - restore it explicitly like the rest of the data rows.
$ws.Range("A13:J20").RowHeight = 20

# ---------------------------------------------------------------------------
# 2. Populate the 8 new config rows.
# ---------------------------------------------------------------------------

# core.test.client.log.request
$ws.Range("A13").Value = "core.test.client.log.request"
$ws.Range("B13").Value = "clientLogRequests"
$ws.Range("C13").Value = "Log requests to be sent from the client (in the client log)"
$ws.Range("G13").Value = "Boolean"
$ws.Range("H13").Value = $true

# core.test.client.log.request.print-fields
$ws.Range("A14").Value = "core.test.client.log.request.print-fields"
$ws.Range("B14").Value = "clientLogRequestsPrintFields"
$ws.Range("C14").Value = "Log each fields of the requests to be sent from the client (in the client log)"
$ws.Range("G14").Value = "Boolean"
$ws.Range("H14").Value = $true

# core.test.client.log.response
$ws.Range("A15").Value = "core.test.client.log.response"
$ws.Range("B15").Value = "clientLogResponses"
$ws.Range("C15").Value = "Log responses received by the client (in the client log)"
$ws.Range("G15").Value = "Boolean"
$ws.Range("H15").Value = $true

# core.test.client.log.response.print-fields
$ws.Range("A16").Value = "core.test.client.log.response.print-fields"
$ws.Range("B16").Value = "clientLogResponsesPrintFields"
$ws.Range("C16").Value = "Log each fields of the responses received by the client (in the client log)"
$ws.Range("G16").Value = "Boolean"
$ws.Range("H16").Value = $true

# core.test.server.log.request
$ws.Range("A17").Value = "core.test.server.log.request"
$ws.Range("B17").Value = "serverLogRequests"
$ws.Range("C17").Value = "Log requests received by the server (in the server log)"
$ws.Range("G17").Value = "Boolean"
$ws.Range("H17").Value = $true

# core.test.server.log.request.print-fields
$ws.Range("A18").Value = "core.test.server.log.request.print-fields"
$ws.Range("B18").Value = "serverLogRequestsPrintFields"
$ws.Range("C18").Value = "Log each fields of the requests received by the server(in the server log)"
$ws.Range("G18").Value = "Boolean"
$ws.Range("H18").Value = $true

# core.test.server.log.response
$ws.Range("A19").Value = "core.test.server.log.response"
$ws.Range("B19").Value = "serverLogResponses"
$ws.Range("C19").Value = "Log responses to be sent from the server (in the server log)"
$ws.Range("G19").Value = "Boolean"
$ws.Range("H19").Value = $true

# core.test.server.log.response.print-fields
$ws.Range("A20").Value = "core.test.server.log.response.print-fields"
$ws.Range("B20").Value = "serverLogResponsesPrintFields"
$ws.Range("C20").Value = "Log each fields of the responses to be sent from the server (in the server log)"
$ws.Range("G20").Value = "Boolean"
$ws.Range("H20").Value = $true

# ---------------------------------------------------------------------------
# 3. Widen columns A/B so the new (longer) names/properties still fit.
#    (ColumnWidth is expressed in characters; target stored widths are
#    34.6640625 / 26.6640625, which after Excel's internal 5-pixel padding
#    correspond to the values below.)
# ---------------------------------------------------------------------------
$ws.Columns("A").ColumnWidth = 33.8333333333333
$ws.Columns("B").ColumnWidth = 25.8333333333333

# ---------------------------------------------------------------------------
# 4. Update the view: scrolled up to the top of the block again, with the
#    newly-added blank-ish row (A46, after renumbering) selected.
# ---------------------------------------------------------------------------
$excel.ActiveWindow.ScrollRow = 7
[void]$ws.Range("A46").Select()
